# Auto-generated edit script: Add data for 2023-11-14
# Applies 172 cell value updates across 48 worksheets of the violent-crime-full-year workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6721
$ws.Range("J3").Value = 7108
$ws.Range("H4").Value = 1710
$ws.Range("J4").Value = 1549
$ws.Range("J5").Value = 559
$ws.Range("J6").Value = 9481
$ws.Range("H7").Value = 26021
$ws.Range("J7").Value = 25418

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J7").Value = 736
$ws.Range("J8").Value = 1595
$ws.Range("J9").Value = 135
$ws.Range("J11").Value = 439
$ws.Range("J15").Value = 303
$ws.Range("J16").Value = 102
$ws.Range("J19").Value = 742
$ws.Range("J20").Value = 532
$ws.Range("J21").Value = 70
$ws.Range("J23").Value = 236
$ws.Range("J24").Value = 81
$ws.Range("J25").Value = 127
$ws.Range("J27").Value = 150
$ws.Range("J29").Value = 1384
$ws.Range("J31").Value = 253
$ws.Range("J33").Value = 1145
$ws.Range("J36").Value = 345
$ws.Range("J37").Value = 781
$ws.Range("J39").Value = 15
$ws.Range("J41").Value = 176
$ws.Range("J42").Value = 1099
$ws.Range("J46").Value = 85
$ws.Range("J51").Value = 311
$ws.Range("J52").Value = 643
$ws.Range("J53").Value = 365
$ws.Range("J54").Value = 486
$ws.Range("J56").Value = 36
$ws.Range("H63").Value = 268
$ws.Range("J63").Value = 82
$ws.Range("J65").Value = 636
$ws.Range("J66").Value = 76
$ws.Range("J71").Value = 84
$ws.Range("J72").Value = 98
$ws.Range("J76").Value = 373
$ws.Range("J82").Value = 34
$ws.Range("J83").Value = 508
$ws.Range("J84").Value = 210
$ws.Range("J85").Value = 1053
$ws.Range("J86").Value = 163
$ws.Range("J88").Value = 263
$ws.Range("J89").Value = 324
$ws.Range("J90").Value = 271
$ws.Range("J91").Value = 293
$ws.Range("J93").Value = 106
$ws.Range("J95").Value = 365
$ws.Range("J96").Value = 277
$ws.Range("J97").Value = 228
$ws.Range("J99").Value = 390
$ws.Range("H101").Value = 26021
$ws.Range("J101").Value = 25418

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J4").Value = 30
$ws.Range("J6").Value = 236
$ws.Range("J7").Value = 736

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 127
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 200
$ws.Range("J7").Value = 439

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 95
$ws.Range("J7").Value = 324

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 285
$ws.Range("J6").Value = 306
$ws.Range("J7").Value = 1053

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 151
$ws.Range("J3").Value = 183
$ws.Range("J7").Value = 643

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 431
$ws.Range("J3").Value = 478
$ws.Range("J6").Value = 561
$ws.Range("J7").Value = 1595

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 190
$ws.Range("J7").Value = 508

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 267
$ws.Range("J3").Value = 380
$ws.Range("J7").Value = 1145

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 130
$ws.Range("J7").Value = 365

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 264
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 781

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 179
$ws.Range("J7").Value = 636

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 390

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 210

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 118
$ws.Range("J6").Value = 229
$ws.Range("J7").Value = 486

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 419
$ws.Range("J3").Value = 490
$ws.Range("J6").Value = 351
$ws.Range("J7").Value = 1384

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 211
$ws.Range("J6").Value = 287
$ws.Range("J7").Value = 742

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 203
$ws.Range("J7").Value = 373

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 233
$ws.Range("J3").Value = 216
$ws.Range("J6").Value = 584
$ws.Range("J7").Value = 1099

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 293

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 181
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 532

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 110
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 345

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 303

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 15

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J3").Value = 43
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 41
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value = 40
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 94
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 83
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("J3").Value = 4
$ws.Range("J6").Value = 34

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 102
